# Sincronização de dados: insere um novo registro de avaliação de garantia
# no topo da série temporal (linha 6), empurrando as linhas seguintes
# (6-28) para baixo (7-29), igual a um novo row puxado da fonte de dados.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insere uma nova linha inteira na posição 6, deslocando as linhas 6..28
# para 7..29 (dimension passa de A1:D28 para A1:D29).
$ws.Range("A6").EntireRow.Insert()

# Preenche a nova linha 6 com o registro mais recente.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = 45986.45500729167
$ws.Range("D6").Value = "MzFmMDgzY2ItODY1MC00N2NkLTg1MjctMzhjMzM4NzRkNDNhOjU3MDE2"
